# Update automatico via Actualizar 07-03-2020 17-06-39
# Appends 7 new rows (543-549) of IGPA index data to the IGPA worksheet,
# extends the IGPA defined name, and moves the selection to reflect the
# new last data row - matching the author's commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IGPA")

# --- Copy existing formatting down onto the new rows before writing values ---
# Rows 543, 544, 548, 549 mirror the "numeric value" look of row 542
# (date in col A / number in col B).
$ws.Range("A542:B542").Copy()
$ws.Range("A543:B544").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A548:B549").PasteSpecial(-4122)   # xlPasteFormats

# Rows 545-547 mirror the date formatting of column A from row 542 ...
$ws.Range("A545:A547").PasteSpecial(-4122)   # xlPasteFormats

# ... and the "--" placeholder look (text, right aligned) of column B,
# taken from an existing "--" cell such as B524.
$ws.Range("B524").Copy()
$ws.Range("B545:B547").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# --- New row values ---
$ws.Range("A543").Value = 44007
$ws.Range("B543").Value = 20306.740000000002

$ws.Range("A544").Value = 44008
$ws.Range("B544").Value = 20100.36

$ws.Range("A545").Value = 44009
$ws.Range("B545").Value = "--"

$ws.Range("A546").Value = 44010
$ws.Range("B546").Value = "--"

$ws.Range("A547").Value = 44011
$ws.Range("B547").Value = "--"

$ws.Range("A548").Value = 44012
$ws.Range("B548").Value = 19969.02

$ws.Range("A549").Value = 44013
$ws.Range("B549").Value = 20315.63

# --- Extend the named range IGPA to cover the new data ---
$wb.Names.Item("IGPA").RefersTo = "=IGPA!`$A`$1:`$B`$549"

# --- Move the active selection to mirror the author's saved state ---
$ws.Range("A551").Select() | Out-Null
